$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SVM)
$ws.Range("B2").Value = 0.672
$ws.Range("D2").Value = 0.67
$ws.Range("F2").Value = 0

# Row 3 (LR)
$ws.Range("B3").Value = 0.672
$ws.Range("D3").Value = 0.6860000000000001

# Row 4 (LDA)
$ws.Range("B4").Value = 0.63
$ws.Range("D4").Value = 0.62

# Row 5 (RF)
$ws.Range("B5").Value = 0.667
$ws.Range("C5").Value = 0.015
$ws.Range("D5").Value = 0.664
$ws.Range("F5").Value = 0.01523745972594797

# Row 6 (AB)
$ws.Range("B6").Value = 0.667
$ws.Range("D6").Value = 0.6860000000000001
$ws.Range("F6").Value = 0

# Row 7 (KNN)
$ws.Range("B7").Value = 0.587
$ws.Range("D7").Value = 0.595

# Row 8 (GNB)
$ws.Range("B8").Value = 0.624
$ws.Range("D8").Value = 0.598
